# Update column C ("Förändrad") for all data rows from 45202 to 45203
# (i.e. bump the date value by one day, from 2023-10-03 to 2023-10-04)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45202) {
        $cell.Value = 45203
    }
}
